$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.254.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.030.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.86"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.50%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.027.33"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.43"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.26%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.249.14"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.534.06"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.61"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +20.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.031.93"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "473.83"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.43"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.84"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.23%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.25%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.64"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.13"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.993"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.16"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +9.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.05"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.60"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "383.23"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.722.02"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.66"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.56"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.14%  "

# Row swaps: Hedera (32) <-> PEPE (33)
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0995"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.27%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.116"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.35%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.84"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.20%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.62"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.51%  "
